# Preparing to write GUI
# Replace the "combined Stats-this session" data (previously 10 players,
# rows 2-11) with a new 4-player data set (rows 2-5), and shrink the
# dependent charts' series ranges from $2:$11 to $2:$5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("combined Stats-this session")
$sheetQ = "'combined Stats-this session'"

# ---------------------------------------------------------------------
# 1. Overwrite rows 2-5 with the new stats.
# ---------------------------------------------------------------------

# Row 2 - Raymond
$ws.Cells.Item(2, 1).Value = "Raymond"
$ws.Cells.Item(2, 2).Value = 60
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = -60
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.587
$ws.Cells.Item(2, 7).Value = 0.331
$ws.Cells.Item(2, 8).Value = 0.0407
$ws.Cells.Item(2, 9).Value = 0.444
$ws.Cells.Item(2, 10).Value = 0.174
$ws.Cells.Item(2, 11).Value = 0.07000000000000001
$ws.Cells.Item(2, 12).Value = 2.08
$ws.Cells.Item(2, 13).Value = 27
$ws.Cells.Item(2, 14).Value = 62
$ws.Cells.Item(2, 15).Value = 109.46
$ws.Cells.Item(2, 16).Value = 93.31999999999999
$ws.Cells.Item(2, 17).Value = 172
# T2 ("Date") is unchanged.

# Row 3 - Cedric
$ws.Cells.Item(3, 1).Value = "Cedric"
$ws.Cells.Item(3, 2).Value = 50
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = -50
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.513
$ws.Cells.Item(3, 7).Value = 0.013
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0.23
$ws.Cells.Item(3, 10).Value = 0.173
$ws.Cells.Item(3, 11).Value = 0.103
$ws.Cells.Item(3, 12).Value = 0.57
$ws.Cells.Item(3, 13).Value = 1
$ws.Cells.Item(3, 14).Value = 2
$ws.Cells.Item(3, 15).Value = 80.69
$ws.Cells.Item(3, 16).Value = 47.72
$ws.Cells.Item(3, 17).Value = 156
$ws.Cells.Item(3, 20).Value = "06/10"

# Row 4 - Fish
$ws.Cells.Item(4, 1).Value = "Fish"
$ws.Cells.Item(4, 2).Value = 40
$ws.Cells.Item(4, 3).Value = 80.84
$ws.Cells.Item(4, 4).Value = 40.84
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.676
$ws.Cells.Item(4, 7).Value = 0.005
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0.161
$ws.Cells.Item(4, 10).Value = 0.198
$ws.Cells.Item(4, 11).Value = 0.077
$ws.Cells.Item(4, 12).Value = 0.28
$ws.Cells.Item(4, 13).Value = 1
$ws.Cells.Item(4, 14).Value = 2
$ws.Cells.Item(4, 15).Value = 212.21
$ws.Cells.Item(4, 16).Value = 85.87
$ws.Cells.Item(4, 17).Value = 182

# Row 5 - Scott
$ws.Cells.Item(5, 1).Value = "Scott"
$ws.Cells.Item(5, 2).Value = 20
$ws.Cells.Item(5, 3).Value = 89.16
$ws.Cells.Item(5, 4).Value = 69.16
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.6860000000000001
$ws.Cells.Item(5, 7).Value = 0.314
$ws.Cells.Item(5, 8).Value = 0.0254
$ws.Cells.Item(5, 9).Value = 0.483
$ws.Cells.Item(5, 10).Value = 0.127
$ws.Cells.Item(5, 11).Value = 0.08500000000000001
$ws.Cells.Item(5, 12).Value = 2.09
$ws.Cells.Item(5, 13).Value = 21
$ws.Cells.Item(5, 14).Value = 37
$ws.Cells.Item(5, 15).Value = 143.51
$ws.Cells.Item(5, 16).Value = 108.51
$ws.Cells.Item(5, 17).Value = 118

# ---------------------------------------------------------------------
# 2. Remove the now-obsolete rows 6-11 (shrinks dimension to A1:T5).
# ---------------------------------------------------------------------
$ws.Rows("6:11").Delete()

# ---------------------------------------------------------------------
# 3. Shrink the dependent chart series from row 11 down to row 5.
#    Only the cat/val ranges move; the series-name (tx) refs are left
#    exactly as they were.
# ---------------------------------------------------------------------
$cos = $ws.ChartObjects()

function Resize-Series($chart, $index, $col) {
    $s = $chart.SeriesCollection().Item($index)
    $s.XValues = "=" + $sheetQ + "!`$A`$2:`$A`$5"
    $s.Values = "=" + $sheetQ + "!`$" + $col + "`$2:`$" + $col + "`$5"
}

# Chart 1: VPIP (F), Pre-flop Raise (G), Three-bet (H)
$chart1 = $cos.Item(1).Chart
Resize-Series $chart1 1 "F"
Resize-Series $chart1 2 "G"
Resize-Series $chart1 3 "H"

# Chart 2: Won at showdown (L)
$chart2 = $cos.Item(2).Chart
Resize-Series $chart2 1 "L"

# Chart 3: C-bets (M), C-bet opportunities (N)
$chart3 = $cos.Item(3).Chart
Resize-Series $chart3 1 "M"
Resize-Series $chart3 2 "N"

# Chart 4: Aggro Frequency (J), Went to showdown (K)
$chart4 = $cos.Item(4).Chart
Resize-Series $chart4 1 "J"
Resize-Series $chart4 2 "K"

# Chart 5: At showdown (O), Before showdown (P)
$chart5 = $cos.Item(5).Chart
Resize-Series $chart5 1 "O"
Resize-Series $chart5 2 "P"

# Chart 6: Hands played (Q)
$chart6 = $cos.Item(6).Chart
Resize-Series $chart6 1 "Q"
